# Swap the data in rows 30 and 31 for the columns that differ between
# the two observation records (A, B, D, E, F, G, H, Q, R, Z, AB).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $cols) {
    $rng30 = $ws.Range($col + "30")
    $rng31 = $ws.Range($col + "31")

    $val30 = $rng30.Value2
    $val31 = $rng31.Value2

    $rng30.Value = $val31
    $rng31.Value = $val30
}
